# fix(import): add antenne column
#
# Adds a new "antenne" column (K) to the measures import template:
#   - K1 header  -> "antenne"
#   - K2 sample  -> "MONTREUIL"
# Leaves the selection on K3 (the cell right below the new sample value),
# matching where a user would land after filling in the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "antenne"
$ws.Range("K2").Value = "MONTREUIL"

$ws.Range("K3").Select()
